{"js": "// Remove the reference to \"SysML\" in the closing sentence of the\n// \"Project adapter.*.ecore\" section:\n//   \"... adapt or extend the SysML concepts which are being exposed by the OSLC adapter.\"\n// becomes\n//   \"... adapt or extend the concepts which are being exposed by the OSLC adapter.\"\n//\n// Note: the document contains a second, unrelated occurrence of \"SysML\"\n// (\"... describe OSLC resources representing SysML concepts\") that must be\n// left untouched, so we search for the unique, longer phrase below rather\n// than the bare word \"SysML\".\n\nconst searchPhrase = \"the SysML concepts which are being exposed by the OSLC adapter\";\nconst replacementPhrase = \"the concepts which are being exposed by the OSLC adapter\";\n\nconst results = context.document.body.search(searchPhrase, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(replacementPhrase, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Remove the reference to \"SysML\" in the closing sentence of the\n# \"Project adapter.*.ecore\" section:\n#   \"... adapt or extend the SysML concepts which are being exposed by the OSLC adapter.\"\n# becomes\n#   \"... adapt or extend the concepts which are being exposed by the OSLC adapter.\"\n#\n# Note: the document contains a second, unrelated occurrence of \"SysML\"\n# (\"... describe OSLC resources representing SysML concepts\") that must be\n# left untouched, so we search/replace the unique, longer phrase below\n# rather than the bare word \"SysML\".\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"the SysML concepts which are being exposed by the OSLC adapter\"\n$find.Replacement.Text = \"the concepts which are being exposed by the OSLC adapter\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
